$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.575558
$ws.Cells.Item(2, 8).Value = 10.726674
$ws.Cells.Item(2, 9).Value = 0.025194653521236
$ws.Cells.Item(2, 10).Value = 0.02519465352123599
$ws.Cells.Item(2, 13).Value = 3.691731666666667
$ws.Cells.Item(2, 14).Value = 11.075195
$ws.Cells.Item(2, 15).Value = 0.04949987503010053
$ws.Cells.Item(2, 16).Value = 0.04949987503010053
$ws.Cells.Item(2, 17).Value = 13.20000069460333
$ws.Cells.Item(2, 18).Value = 118.80000625143
$ws.Cells.Item(2, 19).Value = 0.001247132200727864
$ws.Cells.Item(2, 20).Value = 0.001247132200727864

$ws.Cells.Item(3, 7).Value = 3.575558
$ws.Cells.Item(3, 8).Value = 10.726674
$ws.Cells.Item(3, 9).Value = 0.025194653521236
$ws.Cells.Item(3, 10).Value = 0.02519465352123599
$ws.Cells.Item(3, 15).Value = 0.05430547069958891
$ws.Cells.Item(3, 16).Value = 0.05430547069958892
$ws.Cells.Item(3, 17).Value = 14.481496175888
$ws.Cells.Item(3, 18).Value = 130.333465582992
$ws.Cells.Item(3, 19).Value = 0.001368207518583776
$ws.Cells.Item(3, 20).Value = 0.001368207518583776

$ws.Cells.Item(4, 7).Value = 3.575558
$ws.Cells.Item(4, 8).Value = 10.726674
$ws.Cells.Item(4, 9).Value = 0.025194653521236
$ws.Cells.Item(4, 10).Value = 0.02519465352123599
$ws.Cells.Item(4, 13).Value = 66.27215200000001
$ws.Cells.Item(4, 14).Value = 198.816456
$ws.Cells.Item(4, 15).Value = 0.8885974220704449
$ws.Cells.Item(4, 16).Value = 0.888597422070445
$ws.Cells.Item(4, 17).Value = 236.959923260816
$ws.Cells.Item(4, 18).Value = 2132.639309347344
$ws.Cells.Item(4, 19).Value = 0.02238790416892836
$ws.Cells.Item(4, 20).Value = 0.02238790416892836

$ws.Cells.Item(5, 7).Value = 3.575558
$ws.Cells.Item(5, 8).Value = 10.726674
$ws.Cells.Item(5, 9).Value = 0.025194653521236
$ws.Cells.Item(5, 10).Value = 0.02519465352123599
$ws.Cells.Item(5, 13).Value = 0.5666063333333334
$ws.Cells.Item(5, 14).Value = 1.699819
$ws.Cells.Item(5, 15).Value = 0.007597232199865597
$ws.Cells.Item(5, 16).Value = 0.007597232199865597
$ws.Cells.Item(5, 17).Value = 2.025933808000667
$ws.Cells.Item(5, 18).Value = 18.233404272006
$ws.Cells.Item(5, 19).Value = 0.0001914096329959913
$ws.Cells.Item(5, 20).Value = 0.0001914096329959912

$ws.Cells.Item(6, 9).Value = 0.7460690747908298
$ws.Cells.Item(6, 10).Value = 0.7460690747908298
$ws.Cells.Item(6, 13).Value = 3.691731666666667
$ws.Cells.Item(6, 14).Value = 11.075195
$ws.Cells.Item(6, 15).Value = 0.04949987503010053
$ws.Cells.Item(6, 16).Value = 0.04949987503010053
$ws.Cells.Item(6, 17).Value = 390.8810374058239
$ws.Cells.Item(6, 18).Value = 3517.929336652415
$ws.Cells.Item(6, 19).Value = 0.0369303259659688
$ws.Cells.Item(6, 20).Value = 0.0369303259659688

$ws.Cells.Item(7, 9).Value = 0.7460690747908298
$ws.Cells.Item(7, 10).Value = 0.7460690747908298
$ws.Cells.Item(7, 15).Value = 0.05430547069958891
$ws.Cells.Item(7, 16).Value = 0.05430547069958892
$ws.Cells.Item(7, 19).Value = 0.04051563228092282
$ws.Cells.Item(7, 20).Value = 0.04051563228092282

$ws.Cells.Item(8, 9).Value = 0.7460690747908298
$ws.Cells.Item(8, 10).Value = 0.7460690747908298
$ws.Cells.Item(8, 13).Value = 66.27215200000001
$ws.Cells.Item(8, 14).Value = 198.816456
$ws.Cells.Item(8, 15).Value = 0.8885974220704449
$ws.Cells.Item(8, 16).Value = 0.888597422070445
$ws.Cells.Item(8, 17).Value = 7016.904223774782
$ws.Cells.Item(8, 18).Value = 63152.13801397304
$ws.Cells.Item(8, 19).Value = 0.6629550565456134
$ws.Cells.Item(8, 20).Value = 0.6629550565456135

$ws.Cells.Item(9, 9).Value = 0.7460690747908298
$ws.Cells.Item(9, 10).Value = 0.7460690747908298
$ws.Cells.Item(9, 13).Value = 0.5666063333333334
$ws.Cells.Item(9, 14).Value = 1.699819
$ws.Cells.Item(9, 15).Value = 0.007597232199865597
$ws.Cells.Item(9, 16).Value = 0.007597232199865597
$ws.Cells.Item(9, 17).Value = 59.99235355423812
$ws.Cells.Item(9, 18).Value = 539.9311819881431
$ws.Cells.Item(9, 19).Value = 0.005668059998324827
$ws.Cells.Item(9, 20).Value = 0.005668059998324827

$ws.Cells.Item(10, 7).Value = 32.36130266666667
$ws.Cells.Item(10, 8).Value = 97.08390800000001
$ws.Cells.Item(10, 9).Value = 0.2280292497513723
$ws.Cells.Item(10, 10).Value = 0.2280292497513723
$ws.Cells.Item(10, 13).Value = 3.691731666666667
$ws.Cells.Item(10, 14).Value = 11.075195
$ws.Cells.Item(10, 15).Value = 0.04949987503010053
$ws.Cells.Item(10, 16).Value = 0.04949987503010053
$ws.Cells.Item(10, 17).Value = 119.4692458291178
$ws.Cells.Item(10, 18).Value = 1075.22321246206
$ws.Cells.Item(10, 19).Value = 0.01128741936590051
$ws.Cells.Item(10, 20).Value = 0.01128741936590051

$ws.Cells.Item(11, 7).Value = 32.36130266666667
$ws.Cells.Item(11, 8).Value = 97.08390800000001
$ws.Cells.Item(11, 9).Value = 0.2280292497513723
$ws.Cells.Item(11, 10).Value = 0.2280292497513723
$ws.Cells.Item(11, 15).Value = 0.05430547069958891
$ws.Cells.Item(11, 16).Value = 0.05430547069958892
$ws.Cells.Item(11, 17).Value = 131.0676769371626
$ws.Cells.Item(11, 18).Value = 1179.609092434464
$ws.Cells.Item(11, 19).Value = 0.01238323574102239
$ws.Cells.Item(11, 20).Value = 0.01238323574102239

$ws.Cells.Item(12, 7).Value = 32.36130266666667
$ws.Cells.Item(12, 8).Value = 97.08390800000001
$ws.Cells.Item(12, 9).Value = 0.2280292497513723
$ws.Cells.Item(12, 10).Value = 0.2280292497513723
$ws.Cells.Item(12, 13).Value = 66.27215200000001
$ws.Cells.Item(12, 14).Value = 198.816456
$ws.Cells.Item(12, 15).Value = 0.8885974220704449
$ws.Cells.Item(12, 16).Value = 0.888597422070445
$ws.Cells.Item(12, 17).Value = 2144.653169243339
$ws.Cells.Item(12, 18).Value = 19301.87852319005
$ws.Cells.Item(12, 19).Value = 0.202626203485727
$ws.Cells.Item(12, 20).Value = 0.2026262034857271

$ws.Cells.Item(13, 7).Value = 32.36130266666667
$ws.Cells.Item(13, 8).Value = 97.08390800000001
$ws.Cells.Item(13, 9).Value = 0.2280292497513723
$ws.Cells.Item(13, 10).Value = 0.2280292497513723
$ws.Cells.Item(13, 13).Value = 0.5666063333333334
$ws.Cells.Item(13, 14).Value = 1.699819
$ws.Cells.Item(13, 15).Value = 0.007597232199865597
$ws.Cells.Item(13, 16).Value = 0.007597232199865597
$ws.Cells.Item(13, 17).Value = 18.33611904585023
$ws.Cells.Item(13, 18).Value = 165.025071412652
$ws.Cells.Item(13, 19).Value = 0.00173239115872232
$ws.Cells.Item(13, 20).Value = 0.00173239115872232

$ws.Cells.Item(14, 7).Value = 0.1003386666666667
$ws.Cells.Item(14, 8).Value = 0.301016
$ws.Cells.Item(14, 9).Value = 0.000707021936561918
$ws.Cells.Item(14, 10).Value = 0.0007070219365619179
$ws.Cells.Item(14, 13).Value = 3.691731666666667
$ws.Cells.Item(14, 14).Value = 11.075195
$ws.Cells.Item(14, 15).Value = 0.04949987503010053
$ws.Cells.Item(14, 16).Value = 0.04949987503010053
$ws.Cells.Item(14, 17).Value = 0.3704234331244445
$ws.Cells.Item(14, 18).Value = 3.33381089812
$ws.Cells.Item(14, 19).Value = 0.00003499749750335461
$ws.Cells.Item(14, 20).Value = 0.0000349974975033546

$ws.Cells.Item(15, 7).Value = 0.1003386666666667
$ws.Cells.Item(15, 8).Value = 0.301016
$ws.Cells.Item(15, 9).Value = 0.000707021936561918
$ws.Cells.Item(15, 10).Value = 0.0007070219365619179
$ws.Cells.Item(15, 15).Value = 0.05430547069958891
$ws.Cells.Item(15, 16).Value = 0.05430547069958892
$ws.Cells.Item(15, 17).Value = 0.4063852460586666
$ws.Cells.Item(15, 18).Value = 3.657467214528
$ws.Cells.Item(15, 19).Value = 0.00003839515905992985
$ws.Cells.Item(15, 20).Value = 0.00003839515905992985

$ws.Cells.Item(16, 7).Value = 0.1003386666666667
$ws.Cells.Item(16, 8).Value = 0.301016
$ws.Cells.Item(16, 9).Value = 0.000707021936561918
$ws.Cells.Item(16, 10).Value = 0.0007070219365619179
$ws.Cells.Item(16, 13).Value = 66.27215200000001
$ws.Cells.Item(16, 14).Value = 198.816456
$ws.Cells.Item(16, 15).Value = 0.8885974220704449
$ws.Cells.Item(16, 16).Value = 0.888597422070445
$ws.Cells.Item(16, 17).Value = 6.649659368810668
$ws.Cells.Item(16, 18).Value = 59.846934319296
$ws.Cells.Item(16, 19).Value = 0.000628257870176174
$ws.Cells.Item(16, 20).Value = 0.000628257870176174

$ws.Cells.Item(17, 7).Value = 0.1003386666666667
$ws.Cells.Item(17, 8).Value = 0.301016
$ws.Cells.Item(17, 9).Value = 0.000707021936561918
$ws.Cells.Item(17, 10).Value = 0.0007070219365619179
$ws.Cells.Item(17, 13).Value = 0.5666063333333334
$ws.Cells.Item(17, 14).Value = 1.699819
$ws.Cells.Item(17, 15).Value = 0.007597232199865597
$ws.Cells.Item(17, 16).Value = 0.007597232199865597
$ws.Cells.Item(17, 17).Value = 0.05685252401155557
$ws.Cells.Item(17, 18).Value = 0.5116727161040001
$ws.Cells.Item(17, 19).Value = 0.000005371409822459535
$ws.Cells.Item(17, 20).Value = 0.000005371409822459534
